# B6-PowerPoint.pptx edit — Fri, Jul 03, 2020 3:05:15 PM
#
# 1) Three tables (on the slides that used to carry the custom
#    "Table_0" style {AF65E560-BE96-4D80-96C4-AA6A9B098ED6}) are
#    switched to PowerPoint's built-in "No Style, Table Grid" style
#    {D878C469-E13D-4F08-A18D-4CABB9CE2B34}.
#
# 2) The presentation's theme ("Integral" / Red Violet colours) and the
#    notes-master's theme ("Office Theme" / default colours) are swapped.
#    The two theme parts are identical apart from their colour scheme,
#    so the swap is expressed as applying the Office colour scheme to
#    the (single) editable theme exposed by the object model.

$p = $ppt.ActivePresentation

function HexColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# ---------------------------------------------------------------------
# 1) Retarget the three tables to the built-in table style.
# ---------------------------------------------------------------------
$newTableStyleId = "{D878C469-E13D-4F08-A18D-4CABB9CE2B34}"

for ($idx = 14; $idx -le 16; $idx++) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------
# 2) Swap the theme colours: the deck's theme becomes the plain
#    "Office" colour scheme that used to live on the notes master.
# ---------------------------------------------------------------------
$officeColors = @{
    1  = "000000";  # dk1
    2  = "FFFFFF";  # lt1
    3  = "44546A";  # dk2
    4  = "E7E6E6";  # lt2
    5  = "5B9BD5";  # accent1
    6  = "ED7D31";  # accent2
    7  = "A5A5A5";  # accent3
    8  = "FFC000";  # accent4
    9  = "4472C4";  # accent5
    10 = "70AD47";  # accent6
    11 = "0563C1";  # hlink
    12 = "954F72";  # folHlink
}

$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = HexColor($officeColors[$i])
}
